$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.326.18"
$ws.Range("E2").Value = "  +0.23%  "

$ws.Range("D3").Value = "1.631.39"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'1.002"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "'302.86"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("D7").Value = "'0.3804"
$ws.Range("E7").Value = "  +0.71%  "

$ws.Range("D8").Value = "'51.77"
$ws.Range("E8").Value = "  -0.44%  "

$ws.Range("D9").Value = "'0.3561"
$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("D10").Value = "'0.08134"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").Value = "'22.24"
$ws.Range("E13").Value = "  -1.27%  "

$ws.Range("D14").Value = "'6.416"
$ws.Range("E14").Value = "  -1.67%  "

$ws.Range("D15").Value = "'7.279"
$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").Value = "'0.00001226"
$ws.Range("E16").Value = "  -1.12%  "

$ws.Range("D17").Value = "1.631.68"
$ws.Range("E17").Value = "  +0.67%  "

$ws.Range("D18").Value = "'94.92"
$ws.Range("E18").Value = "  +1.81%  "

$ws.Range("D19").Value = "'0.06944"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("D20").Value = "'6.556"
$ws.Range("E20").Value = "  +2.14%  "

$ws.Range("D21").Value = "'17.28"
$ws.Range("E21").Value = "  -3.05%  "

$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").Value = "'12.41"
$ws.Range("E23").Value = "  -2.24%  "

$ws.Range("D24").Value = "23.390.19"
$ws.Range("E24").Value = "  +0.47%  "

$ws.Range("D25").Value = "'2.540"
$ws.Range("E25").Value = "  +4.12%  "

$ws.Range("D26").Value = "'3.098"
$ws.Range("E26").Value = "  -3.35%  "

$ws.Range("D27").Value = "'21.10"
$ws.Range("E27").Value = "  +0.35%  "

$ws.Range("D28").Value = "'151.72"
$ws.Range("E28").Value = "  +1.83%  "

$ws.Range("D29").Value = "'5.257"
$ws.Range("E29").Value = "  -0.67%  "

$ws.Range("D30").Value = "'133.15"
$ws.Range("E30").Value = "  -0.83%  "

$ws.Range("D31").Value = "1.810.46"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").Value = "'1.079"
$ws.Range("E32").Value = "  +14.50%  "

$ws.Range("D33").Value = "'2.146"
$ws.Range("E33").Value = "  -6.25%  "

$ws.Range("D34").Value = "'6.505"
$ws.Range("E34").Value = "  -3.30%  "

$ws.Range("D35").Value = "'11.44"
$ws.Range("E35").Value = "  +5.31%  "

$ws.Range("D36").Value = "'0.02743"
$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("D37").Value = "'0.2486"
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("D38").Value = "'0.08719"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("D39").Value = "'0.06981"
$ws.Range("E39").Value = "  -1.56%  "

$ws.Range("D40").Value = "'5.922"
$ws.Range("E40").Value = "  -2.37%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6963"
$ws.Range("E41").Value = "  -0.66%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.327"
$ws.Range("E42").Value = "  -2.33%  "

$ws.Range("D43").Value = "'12.15"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").Value = "'15.35"
$ws.Range("E44").Value = "  -4.62%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "'1.002"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6410"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "'2.268"
$ws.Range("E47").Value = "  -1.60%  "

$ws.Range("D48").Value = "'3.952"
$ws.Range("E48").Value = "  -0.68%  "

$ws.Range("D49").Value = "'0.07924"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").Value = "'129.55"
$ws.Range("E50").Value = "  +3.60%  "

$ws.Range("D51").Value = "'1.180"
$ws.Range("E51").Value = "  -1.24%  "

